# Documentatie - Timeline.xlsx: add "1/2 of section 6" progress row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: date, hours, description
# Copy the date-cell number format (style) from B8 down onto B9 first
# (xlPasteFormats = -4122), then fill in the values.
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 45541
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "Terminare jumatate din capitolul 6"

# Move the active selection cursor to C8 (matches diff's sheetView selection)
[void]$ws.Range("C8").Select()
